$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: bill (E2) and phone number (H2) were stored as text, convert to real numbers
$ws.Range("E2").Value = 2400
$ws.Range("H2").Value = 483843992

# New row 3: add order u6745 for Shahroz Ansari
$ws.Range("A3").Value = "u6745"
$ws.Range("B3").Value = "Shahroz Ansari"
$ws.Range("C3").Value = "{'Classic Burger': 1, 'Cheese Burger': 1, 'Chicken Nuggets': 1, 'Onion Rings': 1, 'French Fries': 1, 'Ham Sandwich': 1}"
$ws.Range("D3").Value = "2024-10-06 02:59 PM"
$ws.Range("E3").Value = 2550
$ws.Range("F3").Value = "Bakers street 29"
$ws.Range("G3").Value = "No sauce, extra cheese"
$ws.Range("H3").Value = 3444231978
